# Agent specific features for ordering and delivery added.
# Adds two new columns (N: ordering_period, O: delivery_period) with
# per-agent values for the first three agent rows, and updates the
# current view/selection to reflect the new working area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agent 0 (row 2): ordering_period = 5, delivery_period = 0
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 0

# Agent 1 (row 3): ordering_period = 0, delivery_period = 2
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 2

# Agent 2 (row 4): ordering_period = 0, delivery_period = 3
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 3

# Scroll the view over to show the newly added columns and leave the
# selection on the last entered cell, as would happen after manually
# keying in this data.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("O4").Select()
